$wb = $excel.ActiveWorkbook

# Update "Date studenti" sheet (average from previous year)
$wsDate = $wb.Worksheets.Item("Date studenti")
$wsDate.Range("B2").Value = 7.54
$wsDate.Range("B3").Value = 6.566

# Update "Răspunsuri la formular 1" sheet (Medie column, P)
$wsForm = $wb.Worksheets.Item("Răspunsuri la formular 1")
$wsForm.Range("P2").Value = 7.54
$wsForm.Range("P3").Value = 6.566
